# Update "想去人数" (interested-count) figures on the 展览 / 演出 / 全部类型 sheets
# to reflect the latest scrape (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 309
$ws1.Range("F3").Value = 45
$ws1.Range("F5").Value = 4542
$ws1.Range("F6").Value = 347
$ws1.Range("F8").Value = 284
$ws1.Range("F9").Value = 703
$ws1.Range("F10").Value = 186

# 演出 (Performances) sheet
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 19

# 全部类型 (All types) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 309
$ws4.Range("F3").Value = 45
$ws4.Range("F5").Value = 4542
$ws4.Range("F6").Value = 347
$ws4.Range("F8").Value = 284
$ws4.Range("F9").Value = 703
$ws4.Range("F10").Value = 19
$ws4.Range("F11").Value = 186

$wb.Save()
